$wb = $excel.ActiveWorkbook

# --- Sheet "Runs": insert 4 new rows before current row 6 (Run116), shifting it to row 10 ---
$wsRuns = $wb.Worksheets.Item("Runs")
$wsRuns.Rows.Item(6).Insert()
$wsRuns.Rows.Item(6).Insert()
$wsRuns.Rows.Item(6).Insert()
$wsRuns.Rows.Item(6).Insert()

$wsRuns.Range("A6").Value = 0
$wsRuns.Range("B6").Value = "Run112"
$wsRuns.Range("C6").Value = "0.357 (0.417)"
$wsRuns.Range("D6").Value = "0.561 (0.199)"
$wsRuns.Range("E6").Value = "0.282 (0.450)"
$wsRuns.Range("F6").Value = "0.304 (0.428)"
$wsRuns.Range("G6").Value = "0.720 (0.213)"
$wsRuns.Range("H6").Value = "0.210 (0.408)"
$wsRuns.Range("I6").Value = "0.456 (0.462)"
$wsRuns.Range("J6").Value = "0.593 (0.169)"
$wsRuns.Range("K6").Value = "0.429 (0.496)"
$wsRuns.Range("L6").Value = "0.381 (0.438)"
$wsRuns.Range("M6").Value = "0.594 (0.201)"
$wsRuns.Range("N6").Value = "0.321 (0.467)"

$wsRuns.Range("A7").Value = 0
$wsRuns.Range("B7").Value = "Run113"
$wsRuns.Range("C7").Value = "0.232 (0.368)"
$wsRuns.Range("D7").Value = "0.581 (0.181)"
$wsRuns.Range("E7").Value = "0.143 (0.351)"
$wsRuns.Range("F7").Value = "0.305 (0.390)"
$wsRuns.Range("G7").Value = "0.574 (0.236)"
$wsRuns.Range("H7").Value = "0.180 (0.385)"
$wsRuns.Range("I7").Value = "0.361 (0.410)"
$wsRuns.Range("J7").Value = "0.568 (0.195)"
$wsRuns.Range("K7").Value = "0.272 (0.445)"
$wsRuns.Range("L7").Value = "0.292 (0.392)"
$wsRuns.Range("M7").Value = "0.574 (0.199)"
$wsRuns.Range("N7").Value = "0.193 (0.395)"

$wsRuns.Range("A8").Value = 0
$wsRuns.Range("B8").Value = "Run114"
$wsRuns.Range("C8").Value = "0.456 (0.466)"
$wsRuns.Range("D8").Value = "0.648 (0.117)"
$wsRuns.Range("E8").Value = "0.423 (0.495)"
$wsRuns.Range("F8").Value = "0.345 (0.388)"
$wsRuns.Range("G8").Value = "0.601 (0.237)"
$wsRuns.Range("H8").Value = "0.165 (0.373)"
$wsRuns.Range("I8").Value = "0.588 (0.443)"
$wsRuns.Range("J8").Value = "0.775 (0.135)"
$wsRuns.Range("K8").Value = "0.511 (0.500)"
$wsRuns.Range("L8").Value = "0.499 (0.455)"
$wsRuns.Range("M8").Value = "0.703 (0.177)"
$wsRuns.Range("N8").Value = "0.432 (0.496)"

$wsRuns.Range("A9").Value = 0
$wsRuns.Range("B9").Value = "Run115"
$wsRuns.Range("C9").Value = "0.515 (0.466)"
$wsRuns.Range("D9").Value = "0.620 (0.249)"
$wsRuns.Range("E9").Value = "0.490 (0.500)"
$wsRuns.Range("F9").Value = "0.558 (0.460)"
$wsRuns.Range("G9").Value = "0.769 (0.116)"
$wsRuns.Range("H9").Value = "0.500 (0.501)"
$wsRuns.Range("I9").Value = "0.672 (0.441)"
$wsRuns.Range("J9").Value = "0.529 (0.171)"
$wsRuns.Range("K9").Value = "0.692 (0.462)"
$wsRuns.Range("L9").Value = "0.580 (0.461)"
$wsRuns.Range("M9").Value = "0.630 (0.222)"
$wsRuns.Range("N9").Value = "0.570 (0.495)"

# --- Sheet "Iterations": add new row 2 (Iteration31) ---
$wsIter = $wb.Worksheets.Item("Iterations")
$wsIter.Range("A2").Value = 0
$wsIter.Range("B2").Value = "Iteration31"
$wsIter.Range("C2").Value = "0.39 (0.11)"
$wsIter.Range("D2").Value = "0.60 (0.03)"
$wsIter.Range("E2").Value = "0.33 (0.13)"
$wsIter.Range("F2").Value = "0.38 (0.11)"
$wsIter.Range("G2").Value = "0.67 (0.08)"
$wsIter.Range("H2").Value = "0.26 (0.14)"
$wsIter.Range("I2").Value = "0.52 (0.12)"
$wsIter.Range("J2").Value = "0.62 (0.09)"
$wsIter.Range("K2").Value = "0.48 (0.15)"
$wsIter.Range("L2").Value = "0.44 (0.11)"
$wsIter.Range("M2").Value = "0.63 (0.05)"
$wsIter.Range("N2").Value = "0.38 (0.14)"
